$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.663.28"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.605.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.48%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.51%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.76%  "

# Row 9
$ws.Range("E9").Value = "  +0.19%  "

# Row 10
$ws.Range("E10").Value = "  +0.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.415"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.209.95"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.08%  "

# Row 13
$ws.Range("E13").Value = "  +1.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.12%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.589.78"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.89%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.712.10"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
$ws.Range("E17").Value = "  +0.93%  "

# Row 18
$ws.Range("E18").Value = "  +0.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.07"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.14"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.621"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.94%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.740.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.98%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("E26").Value = "  +2.90%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("D28").ClearFormats()

# Row 29
$ws.Range("E29").Value = "  +0.29%  "

# Row 30
$ws.Range("E30").Value = "  +0.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.599.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.25%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.38%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.54%  "

# Row 34
$ws.Range("E34").Value = "  -2.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.87"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.27%  "

# Row 36
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.72"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.64"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.14"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.25%  "

# Row 40
$ws.Range("E40").Value = "  +0.78%  "

# Row 41
$ws.Range("E41").Value = "  +0.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.900"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("E43").Value = "  -2.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.30%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
$ws.Range("E46").Value = "  -2.14%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.36"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.09%  "

# Row 50
$ws.Range("E50").Value = "  +0.06%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.236"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.10%  "
